$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "extr1".."extr8" rows down by two (from rows 8-15 to
# rows 10-17) to make room for two new rows, carrying data + formatting
# along unchanged.
$ws.Rows("8:15").Copy()
$ws.Rows("10:17").PasteSpecial(-4104)
$excel.CutCopyMode = $false

# Fill the freed rows 8 and 9 with the new "line7"/"line8" data (the
# cell formatting there is already correct, inherited from the original
# sheet, so plain value assignment is enough).
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $false

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# Renumber column A (a simple 0-based row index) for the shifted rows,
# and fix up the two "in_service" flags that flipped on the shifted
# "extr1"/"extr2" rows (now rows 10 and 11).
$ws.Range("A10").Value = 8
$ws.Range("E10").Value = $true

$ws.Range("A11").Value = 9
$ws.Range("E11").Value = $true

$ws.Range("A12").Value = 10
$ws.Range("A13").Value = 11
$ws.Range("A14").Value = 12
$ws.Range("A15").Value = 13
$ws.Range("A16").Value = 14
$ws.Range("A17").Value = 15

# Rows 16/17 are brand new (past the original A1:E15 extent), so give
# column A there the same bold/border style used by the rest of column A.
$ws.Range("A2").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
